# Updates the "cryptos" price/volume snapshot (Price = column D,
# Volume(1h) = column E) to the latest scrape values. A handful of coins
# also swapped ranking position with their neighbour (rows 30/31, 40/41,
# 43/44), so their Coin name (B) and Link (C) cells are rewritten too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("D2").Value = "62.282.21"
    $ws.Range("E2").Value = "  +0.00%  "
    $ws.Range("D3").Value = "3.010.26"
    $ws.Range("E3").Value = "  +0.12%  "
    $ws.Range("E4").Value = "  +0.05%  "
    $ws.Range("D5").Value = "593.04"
    $ws.Range("E5").Value = "  +1.89%  "
    $ws.Range("D6").Value = "147.16"
    $ws.Range("E6").Value = "  +1.06%  "
    $ws.Range("E7").Value = "  +0.09%  "
    $ws.Range("D8").Value = "3.006.60"
    $ws.Range("E8").Value = "  -0.02%  "
    $ws.Range("D9").Value = "0.516"
    $ws.Range("E9").Value = "  -2.11%  "
    $ws.Range("D10").Value = "6.25"
    $ws.Range("E10").Value = "  +8.18%  "
    $ws.Range("E11").Value = "  +0.34%  "
    $ws.Range("D12").Value = "0.456"
    $ws.Range("E12").Value = "  -1.61%  "
    $ws.Range("D13").Value = "0.0000231"
    $ws.Range("E13").Value = "  +1.46%  "
    $ws.Range("D14").Value = "34.38"
    $ws.Range("E14").Value = "  -0.21%  "
    $ws.Range("E15").Value = "  +2.41%  "
    $ws.Range("D16").Value = "3.510.06"
    $ws.Range("E16").Value = "  +0.34%  "
    $ws.Range("D17").Value = "62.243.14"
    $ws.Range("E17").Value = "  -0.05%  "
    $ws.Range("D18").Value = "6.98"
    $ws.Range("E18").Value = "  -1.85%  "
    $ws.Range("D19").Value = "3.012.52"
    $ws.Range("E19").Value = "  +0.16%  "
    $ws.Range("D20").Value = "449.42"
    $ws.Range("E20").Value = "  -1.19%  "
    $ws.Range("D21").Value = "14.09"
    $ws.Range("E21").Value = "  +0.96%  "
    $ws.Range("D22").Value = "0.685"
    $ws.Range("E22").Value = "  -0.36%  "
    $ws.Range("E23").Value = "  -0.61%  "
    $ws.Range("D24").Value = "81.95"
    $ws.Range("E24").Value = "  +0.42%  "
    $ws.Range("D25").Value = "11.00"
    $ws.Range("E25").Value = "  +10.00%  "
    $ws.Range("D26").Value = "2.23"
    $ws.Range("E26").Value = "  +0.63%  "
    $ws.Range("D27").Value = "12.11"
    $ws.Range("E27").Value = "  -1.94%  "
    $ws.Range("E28").Value = "  -0.03%  "
    $ws.Range("D29").Value = "2.69"
    $ws.Range("E29").Value = "  +3.42%  "
    $ws.Range("B30").Value = "FirstDigitalUSD"
    $ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
    $ws.Range("D30").Value = "1.00"
    $ws.Range("E30").Value = "  +0.07%  "
    $ws.Range("B31").Value = "NEARProtocol"
    $ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
    $ws.Range("D31").Value = "7.27"
    $ws.Range("E31").Value = "  +4.30%  "
    $ws.Range("E32").Value = "  -0.03%  "
    $ws.Range("D33").Value = "27.42"
    $ws.Range("E33").Value = "  -3.16%  "
    $ws.Range("E34").Value = "  +1.63%  "
    $ws.Range("D35").Value = "0.0₃0845"
    $ws.Range("E35").Value = "  +5.88%  "
    $ws.Range("E36").Value = "  -0.36%  "
    $ws.Range("D37").Value = "5.81"
    $ws.Range("E37").Value = "  +0.86%  "
    $ws.Range("D38").Value = "50.24"
    $ws.Range("E38").Value = "  -0.02%  "
    $ws.Range("D39").Value = "2.05"
    $ws.Range("E39").Value = "  -3.12%  "
    $ws.Range("B40").Value = "Cosmos"
    $ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
    $ws.Range("D40").Value = "8.96"
    $ws.Range("E40").Value = "  -2.38%  "
    $ws.Range("B41").Value = "dogwifhat"
    $ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
    $ws.Range("D41").Value = "2.95"
    $ws.Range("E41").Value = "  +2.02%  "
    $ws.Range("E42").Value = "  +5.74%  "
    $ws.Range("B43").Value = "Arweave"
    $ws.Range("C43").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
    $ws.Range("D43").Value = "40.95"
    $ws.Range("E43").Value = "  +10.37%  "
    $ws.Range("B44").Value = "Bittensor"
    $ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
    $ws.Range("D44").Value = "396.89"
    $ws.Range("E44").Value = "  +1.59%  "
    $ws.Range("D45").Value = "0.279"
    $ws.Range("E45").Value = "  +3.81%  "
    $ws.Range("D46").Value = "0.0352"
    $ws.Range("E46").Value = "  -1.35%  "
    $ws.Range("D47").Value = "2.709.28"
    $ws.Range("E47").Value = "  -0.40%  "
    $ws.Range("D48").Value = "132.38"
    $ws.Range("E48").Value = "  +3.07%  "
    $ws.Range("D50").Value = "2.18"
    $ws.Range("E50").Value = "  -0.65%  "
    $ws.Range("E51").Value = "  -1.62%  "
